$d = $word.ActiveDocument

$replacements = @(
    @{old = "456×7=3192"; new = "694×3=2082"},
    @{old = "381×9=3429"; new = "965×7=6755"},
    @{old = "635×8=5080"; new = "119×9=1071"},
    @{old = "247×9=2223"; new = "582×5=2910"},
    @{old = "291×9=2619"; new = "913×2=1826"},
    @{old = "265×4=1060"; new = "259×4=1036"},
    @{old = "644×7=4508"; new = "683×5=3415"},
    @{old = "841×4=3364"; new = "383×3=1149"},
    @{old = "196×3=588";  new = "851×7=5957"},
    @{old = "314×3=942";  new = "809×6=4854"},
    @{old = "248×3=744";  new = "966×2=1932"},
    @{old = "126×4=504";  new = "917×4=3668"},
    @{old = "207×5=1035"; new = "379×6=2274"},
    @{old = "225×4=900";  new = "701×9=6309"},
    @{old = "178×5=890";  new = "258×6=1548"},
    @{old = "429×6=2574"; new = "773×3=2319"},
    @{old = "869×2=1738"; new = "362×6=2172"},
    @{old = "629×4=2516"; new = "616×3=1848"},
    @{old = "654×8=5232"; new = "563×7=3941"},
    @{old = "124×7=868";  new = "676×8=5408"},
    @{old = "620×2=1240"; new = "906×9=8154"},
    @{old = "744×9=6696"; new = "977×7=6839"},
    @{old = "886×3=2658"; new = "430×2=860"},
    @{old = "813×2=1626"; new = "609×5=3045"},
    @{old = "108×2=216";  new = "416×2=832"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
